# Apply the committed changes to Test.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1
$ws2 = $wb.Worksheets.Item(2)   # Sheet2

# --- Sheet1: populate new sample rows under the existing header row ---
$ws1.Range("A2").Value = "a"
$ws1.Range("A3").Value = "b"
$ws1.Range("A4").Value = "c"

# --- Sheet2: replace the long placeholder text with short sample values
#     and add the matching new rows ---
$ws2.Range("A1").Value = "a"
$ws2.Range("A2").Value = "b"
$ws2.Range("A3").Value = "c"
$ws2.Range("A4").Value = "d"

# Clear out the old text that used to live in A5/D5 (keep their formatting)
$ws2.Range("A5").Value = ""
$ws2.Range("D5").Value = ""

# Row 8 gained the same wrap-text style as column A elsewhere, with an
# explicit 16pt row height
$ws2.Range("A8").WrapText = $true
$ws2.Rows.Item(8).RowHeight = 16

# --- Update the saved selections / active sheet ---
# Sheet2 keeps a simple selection at B14 but is no longer the active tab
$ws2.Range("B14").Select()

# Sheet1 becomes the active tab, selection moves to C6
$ws1.Activate()
$ws1.Range("C6").Select()
